# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2404
#   *_new  -> *_FV2410
# Also turn the data range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) -------------------------------------------------
$headerRenames = @{
    "A1" = "Segmentname_FV2404";
    "B1" = "Segmentgruppe_FV2404";
    "C1" = "Segment_FV2404";
    "D1" = "Datenelement_FV2404";
    "E1" = "Segment ID_FV2404";
    "F1" = "Code_FV2404";
    "G1" = "Qualifier_FV2404";
    "H1" = "Beschreibung_FV2404";
    "I1" = "Bedingungsausdruck_FV2404";
    "J1" = "Bedingung_FV2404";
    "L1" = "Segmentname_FV2410";
    "M1" = "Segmentgruppe_FV2410";
    "N1" = "Segment_FV2410";
    "O1" = "Datenelement_FV2410";
    "P1" = "Segment ID_FV2410";
    "Q1" = "Code_FV2410";
    "R1" = "Qualifier_FV2410";
    "S1" = "Beschreibung_FV2410";
    "T1" = "Bedingungsausdruck_FV2410";
    "U1" = "Bedingung_FV2410";
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# --- 2. Turn A1:U84 into an Excel Table (adds xl/tables/table1.xml + tableParts) -------
$dataRange = $ws.Range("A1:U84")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row (row 1) ---------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
